# Fruta / hortaliza, semanal
# Update the weekly Higo (fig) price records for Vega Central Mapocho de
# Santiago: dates, volumes, prices and origin labels are refreshed with the
# latest reported figures for the period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Primera)
$ws.Range("D2").Value = 44300
$ws.Range("M2").Value = 100
$ws.Range("R2").Value = "Región Metropolitana"

# Row 3 (Segunda)
$ws.Range("D3").Value = 44300
$ws.Range("M3").Value = 80
$ws.Range("R3").Value = "Región Metropolitana"

# Row 4 (Primera)
$ws.Range("D4").Value = 44299
$ws.Range("M4").Value = 80
$ws.Range("N4").Value = 15000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 15000
$ws.Range("R4").Value = "Provincia de Santiago"
$ws.Range("S4").Value = 2143

# Row 5 (Segunda)
$ws.Range("D5").Value = 44299
$ws.Range("M5").Value = 75
$ws.Range("R5").Value = "Provincia de Santiago"

# Row 6 (Primera)
$ws.Range("D6").Value = 44322
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 12000
$ws.Range("P6").Value = 12000
$ws.Range("S6").Value = 1714

# Row 7 (Segunda)
$ws.Range("D7").Value = 44322
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 8000
$ws.Range("O7").Value = 8000
$ws.Range("P7").Value = 8000
$ws.Range("S7").Value = 1143

# Row 8 (Primera)
$ws.Range("D8").Value = 44320
$ws.Range("M8").Value = 20
$ws.Range("N8").Value = 12000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 12000
$ws.Range("S8").Value = 1714

# Row 9 (Segunda)
$ws.Range("D9").Value = 44320
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 1143

# Row 12 (Primera)
$ws.Range("D12").Value = 44292
$ws.Range("M12").Value = 25
$ws.Range("N12").Value = 16000
$ws.Range("O12").Value = 16000
$ws.Range("P12").Value = 16000
$ws.Range("S12").Value = 2286

# Row 13 (Segunda)
$ws.Range("D13").Value = 44292
$ws.Range("M13").Value = 30
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("S13").Value = 2143

# Row 14 (Primera)
$ws.Range("D14").Value = 44301
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 14000
$ws.Range("O14").Value = 14000
$ws.Range("P14").Value = 14000
$ws.Range("S14").Value = 2000

# Row 15 (Segunda)
$ws.Range("D15").Value = 44301
$ws.Range("M15").Value = 80
$ws.Range("N15").Value = 12000
$ws.Range("O15").Value = 12000
$ws.Range("P15").Value = 12000
$ws.Range("S15").Value = 1714
